$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.318.24"

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.832.70"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.42%  "

# Row 5 - BNB
$ws.Range("D5").Value = "235.75"
$ws.Range("E5").Value = "  -1.43%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.6040"
$ws.Range("E6").Value = "  -2.95%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.37%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -4.54%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.2777"
$ws.Range("E9").Value = "  -3.50%  "

# Row 10 - Solana
$ws.Range("D10").Value = "23.67"
$ws.Range("E10").Value = "  -3.96%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07622"
$ws.Range("E11").Value = "  -1.32%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.850.32"
$ws.Range("E12").Value = "  +0.85%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "4.766"
$ws.Range("E13").Value = "  -3.44%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.6332"
$ws.Range("E14").Value = "  -4.03%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.000009894"
$ws.Range("E15").Value = "  -4.48%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "78.05"

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.032.34"
$ws.Range("E17").Value = "  -0.97%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "5.613"
$ws.Range("E18").Value = "  -9.85%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "218.56"
$ws.Range("E19").Value = "  -7.49%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.37%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -4.56%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "6.936"
$ws.Range("E22").Value = "  -3.69%  "

# Row 23 - BinanceUSD
$ws.Range("E23").Value = "  -0.01%  "

# Row 24 - Monero
$ws.Range("D24").Value = "156.49"
$ws.Range("E24").Value = "  -0.46%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "8.008"
$ws.Range("E25").Value = "  -4.75%  "

# Row 26 - Stellar
$ws.Range("D26").Value = "0.1296"
$ws.Range("E26").Value = "  -2.47%  "

# Row 27 - EthereumClassic: no change

# Row 28 - Hedera
$ws.Range("D28").Value = "0.06454"
$ws.Range("E28").Value = "  -5.83%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "1.427"
$ws.Range("E29").Value = "  -3.23%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "1.446"
$ws.Range("E30").Value = "  -2.19%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "3.852"
$ws.Range("E31").Value = "  -2.05%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "3.811"
$ws.Range("E32").Value = "  -4.98%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "1.735"
$ws.Range("E33").Value = "  -0.34%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "1.098"
$ws.Range("E34").Value = "  -4.72%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "0.6516"
$ws.Range("E35").Value = "  -3.81%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.50%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "2.757"
$ws.Range("E37").Value = "  -0.96%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01759"
$ws.Range("E38").Value = "  -3.27%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "6.619"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40 - Maker
$ws.Range("D40").Value = "1.147.35"
$ws.Range("E40").Value = "  -6.92%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.8961"
$ws.Range("E41").Value = "  -4.90%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.23%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.998.49"

# Row 44 - Quant
$ws.Range("D44").Value = "101.02"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45 - Aave
$ws.Range("D45").Value = "62.44"
$ws.Range("E45").Value = "  -4.02%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -3.38%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "1.624"
$ws.Range("E47").Value = "  -3.62%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "8.545"
$ws.Range("E48").Value = "  -2.93%  "

# Row 49 - now Aptos (was Mantle)
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "6.453"
$ws.Range("E49").Value = "  -6.11%  "

# Row 50 - now Mantle (was Aptos)
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.4548"
$ws.Range("E50").Value = "  -0.44%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.05500"
$ws.Range("E51").Value = "  -2.39%  "
